# Remove the duplicated "고양이 급식기/급수기" category row (row 38).
# Deleting the entire row shifts rows 39-42 up by one, and Excel
# automatically drops the now-unused shared string and recompacts the
# shared-strings table (count/uniqueCount 44 -> 43).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(38).Delete()

# Update the saved view state: the sheet had scrolled to show row 16 at
# the top with H40 selected; after the edit the sheet is scrolled back
# to the top (no topLeftCell override) and H30 is selected.
$ws.Application.Goto($ws.Range("A1"))
[void]$ws.Range("H30").Select()
